$d = $word.ActiveDocument

# Step 1: change the lone space run (with <w:cs/>) before the ellipsis run to "ที่ "
$d.Content.Find.Execute("ตาราง … Activity", $false, $false, $false, $false, $false, $true, 1, $false, "ตารางที่ 1 Activity", 2)
